# Fruta / hortaliza, semanal
# Insert a new weekly record at row 101 (pushing existing rows 101-124
# down to 102-125) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 101; this shifts rows 101-124
# down to 102-125 and extends the sheet dimension to A1:R125.
$ws.Rows.Item(101).Insert()

# Populate the newly inserted row 101 with the new record's data.
$ws.Range("A101").Value2 = 5
$ws.Range("B101").Value2 = "Macroferia Regional de Talca"
$ws.Range("C101").Value2 = "Maule"
$ws.Range("D101").Value2 = 44463
$ws.Range("E101").Value2 = 7
$ws.Range("F101").Value2 = 100112021
$ws.Range("G101").Value2 = "Ají"
$ws.Range("H101").Value2 = "Americana (o)"
$ws.Range("I101").Value2 = "Primera"
$ws.Range("J101").Value2 = 100
$ws.Range("K101").Value2 = 80000
$ws.Range("L101").Value2 = 80000
$ws.Range("M101").Value2 = 80000
$ws.Range("N101").Value2 = "$/caja 25 kilos"
$ws.Range("O101").Value2 = "Provincia del Elquí"
$ws.Range("P101").Value2 = 3200
$ws.Range("Q101").Value2 = 25
$ws.Range("R101").Value2 = "Hortaliza"
